$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing ACTION / ACTION2 columns (E:F),
# pushing them to G:H and leaving new blank columns at E:F.
$ws.Columns("E:F").Insert()

# New condition-column headers (row 7) - reuse the "CONDITION" label (no numbering)
$ws.Range("E7").Value = "CONDITION"
$ws.Range("F7").Value = "CONDITION"

# New condition-column templates (row 8)
$ws.Range("E8").Value = "customer.getBalance() >= 1000"
$ws.Range("F8").Value = "customer.getCreditScore() >= 700"

# Update the YoungAdultRule row's ACTION value (now in column G) to the new test value
$ws.Range("G9").Value = "44sds"
